$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update shared G/H values across existing rows 2-4 (ligand average/total expression values)
$ws.Range("G2:G4").Value = 0.135749
$ws.Range("H2:H4").Value = 0.407247

# Row 2 specific updates
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.7144740000000001
$ws.Range("N2").Value = 2.143422
$ws.Range("O2").Value = 0.138796410342318
$ws.Range("P2").Value = 0.138796410342318
$ws.Range("Q2").Value = 0.09698913102600001
$ws.Range("R2").Value = 0.8729021792340002
$ws.Range("S2").Value = 0.138796410342318
$ws.Range("T2").Value = 0.138796410342318

# Row 3 specific updates
$ws.Range("O3").Value = 0.8044215857867821
$ws.Range("P3").Value = 0.8044215857867821
$ws.Range("Q3").Value = 0.5621193688770001
$ws.Range("R3").Value = 5.059074319893001
$ws.Range("S3").Value = 0.8044215857867821
$ws.Range("T3").Value = 0.8044215857867821

# Row 4 specific updates
$ws.Range("M4").Value = 0.2847646666666667
$ws.Range("N4").Value = 0.8542940000000001
$ws.Range("O4").Value = 0.05531945672713084
$ws.Range("P4").Value = 0.05531945672713083
$ws.Range("Q4").Value = 0.03865651873533334
$ws.Range("R4").Value = 0.3479086686180001
$ws.Range("S4").Value = 0.05531945672713084
$ws.Range("T4").Value = 0.05531945672713083

# New row 5
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Fgf5"
$ws.Range("C5").Value = "Fgfr2"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.135749
$ws.Range("H5").Value = 0.407247
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.007528666666666667
$ws.Range("N5").Value = 0.022586
$ws.Range("O5").Value = 0.00146254714376898
$ws.Range("P5").Value = 0.00146254714376898
$ws.Range("Q5").Value = 0.001022008971333333
$ws.Range("R5").Value = 0.009198080742000002
$ws.Range("S5").Value = 0.00146254714376898
$ws.Range("T5").Value = 0.00146254714376898
